# Added periodic & upfront related scenarios
# The "repaymentstrategy" input on the ProductLoanInput sheet (cell B17)
# is changed from "Mifos style" to a new option describing the
# periodic/upfront order: "Penalties, Fees, Interest, Principal order".
# The cell also picks up a left/top-aligned style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

$cell = $ws.Range("B17")
$cell.Value = "Penalties, Fees, Interest, Principal order"
$cell.HorizontalAlignment = -4131  # xlLeft
$cell.VerticalAlignment = -4160    # xlTop

$cell.Select()
